# Applies the "Natmi following Dr Hou advice" edit:
# Expands the Cxcl12-Cd4 sending/target cluster grid from 3 rows to a full
# 3x3 ECs/FAPs/sCs x ECs/FAPs/sCs grid (9 data rows) with recomputed metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    , @("ECs", "Cxcl12", "Cd4", "ECs", 3, 1, 109.026058, 327.078174, 0.3049840938689738, 0.3049840938689738, 2, 0.6666666666666666, 0.539832, 1.619496, 0.3224687287194008, 0.3224687287194007, 58.85575494225599, 529.701794480304, 0.09834783302956637, 0.09834783302956636)
    , @("ECs", "Cxcl12", "Cd4", "FAPs", 3, 1, 109.026058, 327.078174, 0.3049840938689738, 0.3049840938689738, 3, 1, 1.034252666666666, 3.102758, 0.6178109904463799, 0.6178109904463799, 112.7604912226546, 1014.844421003892, 0.1884225251035824, 0.1884225251035824)
    , @("ECs", "Cxcl12", "Cd4", "sCs", 3, 1, 109.026058, 327.078174, 0.3049840938689738, 0.3049840938689738, 2, 0.6666666666666666, 0.09997533333333335, 0.299926, 0.05972028083421942, 0.05972028083421941, 10.89991649056933, 98.09924841512401, 0.01821373573582505, 0.01821373573582505)
    , @("FAPs", "Cxcl12", "Cd4", "ECs", 3, 1, 89.97721833333333, 269.931655, 0.2516978134001918, 0.2516978134001917, 2, 0.6666666666666666, 0.539832, 1.619496, 0.3224687287194008, 0.3224687287194007, 48.57258172731999, 437.15323554588, 0.0811646739086128, 0.08116467390861277)
    , @("FAPs", "Cxcl12", "Cd4", "FAPs", 3, 1, 89.97721833333333, 269.931655, 0.2516978134001918, 0.2516978134001917, 3, 1, 1.034252666666666, 3.102758, 0.6178109904463799, 0.6178109904463799, 93.05917800049886, 837.5326020044898, 0.1555016753899606, 0.1555016753899605)
    , @("FAPs", "Cxcl12", "Cd4", "sCs", 3, 1, 89.97721833333333, 269.931655, 0.2516978134001918, 0.2516978134001917, 2, 0.6666666666666666, 0.09997533333333335, 0.299926, 0.05972028083421942, 0.05972028083421941, 8.995502395281111, 80.95952155753, 0.01503146410161841, 0.0150314641016184)
    , @("sCs", "Cxcl12", "Cd4", "ECs", 3, 1, 158.477852, 475.433556, 0.4433180927308344, 0.4433180927308344, 2, 0.6666666666666666, 0.539832, 1.619496, 0.3224687287194008, 0.3224687287194007, 85.55141580086399, 769.9627422077759, 0.1429562217812216, 0.1429562217812216)
    , @("sCs", "Cxcl12", "Cd4", "FAPs", 3, 1, 158.477852, 475.433556, 0.4433180927308344, 0.4433180927308344, 3, 1, 1.034252666666666, 3.102758, 0.6178109904463799, 0.6178109904463799, 163.9061410386053, 1475.155269347448, 0.2738867899528369, 0.2738867899528369)
    , @("sCs", "Cxcl12", "Cd4", "sCs", 3, 1, 158.477852, 475.433556, 0.4433180927308344, 0.4433180927308344, 2, 0.6666666666666666, 0.09997533333333335, 0.299926, 0.05972028083421942, 0.05972028083421941, 15.84387607965067, 142.594884716856, 0.02647508099677596, 0.02647508099677595)
)

$rowIdx = 2
foreach ($r in $data) {
    for ($col = 0; $col -lt $r.Length; $col++) {
        $ws.Cells.Item($rowIdx, $col + 1).Value = $r[$col]
    }
    $rowIdx++
}
